$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.170084597650316
$ws.Range("D2").Value = 0.01388653062972622
$ws.Range("E2").Value = 0.8165701030002879
$ws.Range("F2").Value = 0.3243570632220667
$ws.Range("G2").Value = 0.1905593167315516
$ws.Range("H2").Value = 0.3304723969503129
$ws.Range("I2").Value = 0.9194504041202691
$ws.Range("L2").Value = 0.2966425630251877
$ws.Range("M2").Value = 0.3067684467381184
$ws.Range("O2").Value = 0.947654812267217
$ws.Range("B3").Value = 1.06926235010269
$ws.Range("D3").Value = 0.01210636417523148
$ws.Range("E3").Value = 0.7790820133023999
$ws.Range("F3").Value = 0.3161381378764005
$ws.Range("G3").Value = 0.183244661425519
$ws.Range("H3").Value = 0.3309351845615112
$ws.Range("I3").Value = 0.9128897523368948
$ws.Range("L3").Value = 0.2654922027279838
$ws.Range("M3").Value = 0.2780572639620829
$ws.Range("O3").Value = 0.9328345901014359
$ws.Range("B4").Value = 1.007170064841716
$ws.Range("D4").Value = 0.01100983762897556
$ws.Range("E4").Value = 0.75626965982093
$ws.Range("F4").Value = 0.311398351808144
$ws.Range("G4").Value = 0.1789781534719168
$ws.Range("H4").Value = 0.3314581068826001
$ws.Range("I4").Value = 0.9095999869402647
$ws.Range("L4").Value = 0.2463105021988241
$ws.Range("M4").Value = 0.2603769395444502
$ws.Range("O4").Value = 0.9247010386883403
$ws.Range("B5").Value = 0.9818216907802935
$ws.Range("D5").Value = 0.01056215154336115
$ws.Range("E5").Value = 0.7470267848166543
$ws.Range("F5").Value = 0.3095437240104175
$ws.Range("G5").Value = 0.1772956779219967
$ws.Range("H5").Value = 0.3317311644962118
$ws.Range("I5").Value = 0.9084455961756461
$ws.Range("L5").Value = 0.2384804914708809
$ws.Range("M5").Value = 0.2531596271328169
$ws.Range("O5").Value = 0.9216285858528721
$ws.Range("B6").Value = 0.9776099311580424
$ws.Range("D6").Value = 0.01048776380162053
$ws.Range("E6").Value = 0.7454952837722288
$ws.Range("F6").Value = 0.3092404003733478
$ws.Range("G6").Value = 0.1770196859382409
$ws.Range("H6").Value = 0.3317801248044603
$ws.Range("I6").Value = 0.9082651727985009
$ws.Range("L6").Value = 0.2371795373902899
$ws.Range("M6").Value = 0.251960462311807
$ws.Range("O6").Value = 0.921132999144092
$ws.Range("B7").Value = 1.006828388723306
$ws.Range("D7").Value = 0.01100380334106177
$ws.Range("E7").Value = 0.7561447889595314
$ws.Range("F7").Value = 0.3113730287116425
$ws.Range("G7").Value = 0.1789552360851516
$ws.Range("H7").Value = 0.3314615467518962
$ws.Range("I7").Value = 0.9095836637530041
$ws.Range("L7").Value = 0.2462049570596463
$ws.Range("M7").Value = 0.2602796540019838
$ws.Range("O7").Value = 0.9246586237271828
$ws.Range("B8").Value = 1.135360921866948
$ws.Range("D8").Value = 0.01327347655286815
$ws.Range("E8").Value = 0.8036025474336554
$ws.Range("F8").Value = 0.3214593814507438
$ws.Range("G8").Value = 0.187990372518513
$ws.Range("H8").Value = 0.3305823620832342
$ws.Range("I8").Value = 0.9170353126434136
$ws.Range("L8").Value = 0.2859137242949998
$ws.Range("M8").Value = 0.2968798321252706
$ws.Range("O8").Value = 0.9423437306422215
$ws.Range("B9").Value = 1.385864484399463
$ws.Range("D9").Value = 0.01769514543522632
$ws.Range("E9").Value = 0.8982294911903921
$ws.Range("F9").Value = 0.3436842488750287
$ws.Range("G9").Value = 0.2075079715846613
$ws.Range("H9").Value = 0.330756745311561
$ws.Range("I9").Value = 0.9374884013926277
$ws.Range("L9").Value = 0.3633224689406518
$ws.Range("M9").Value = 0.3682243614249359
$ws.Range("O9").Value = 0.9847338378379504
$ws.Range("B10").Value = 1.568898085643696
$ws.Range("D10").Value = 0.02092439368308874
$ws.Range("E10").Value = 0.9686210568539053
$ws.Range("F10").Value = 0.3615228318759875
$ws.Range("G10").Value = 0.2229692771830969
$ws.Range("H10").Value = 0.3320482836003436
$ws.Range("I10").Value = 0.9560547370952861
$ws.Range("L10").Value = 0.4198909809427676
$ws.Range("M10").Value = 0.4203594844888627
$ws.Range("O10").Value = 1.020643079153331
$ws.Range("B11").Value = 1.651932464596882
$ws.Range("D11").Value = 0.02238894557572735
$ws.Range("E11").Value = 1.000815774281065
$ws.Range("F11").Value = 0.3699701231465866
$ws.Range("G11").Value = 0.230251914591193
$ws.Range("H11").Value = 0.3328897953833376
$ws.Range("I11").Value = 0.9652657063087844
$ws.Range("L11").Value = 0.4455550416425069
$ws.Range("M11").Value = 0.4440119986782634
$ws.Range("O11").Value = 1.038028000244935
$ws.Range("B12").Value = 1.683341129237022
$ws.Range("D12").Value = 0.02294286034246795
$ws.Range("E12").Value = 1.013030439051633
$ws.Range("F12").Value = 0.3732169801845799
$ws.Range("G12").Value = 0.2330458836958371
$ws.Range("H12").Value = 0.3332450711589701
$ws.Range("I12").Value = 0.9688632959580588
$ws.Range("L12").Value = 0.4552629080317843
$ws.Range("M12").Value = 0.4529589509839411
$ws.Range("O12").Value = 1.04476317877436
$ws.Range("B13").Value = 1.676578280752665
$ws.Range("D13").Value = 0.02282359572332382
$ws.Range("E13").Value = 1.010398786247293
$ws.Range("F13").Value = 0.3725155704238148
$ws.Range("G13").Value = 0.232442537159784
$ws.Range("H13").Value = 0.3331669263309749
$ws.Range("I13").Value = 0.9680836228621388
$ws.Range("L13").Value = 0.4531726242863954
$ws.Range("M13").Value = 0.451032506442445
$ws.Range("O13").Value = 1.043305867342241
$ws.Range("B14").Value = 1.65451717899532
$ws.Range("D14").Value = 0.022434530315671
$ws.Range("E14").Value = 1.00182022626764
$ws.Range("F14").Value = 0.3702362796029774
$ws.Range("G14").Value = 0.2304810487348163
$ws.Range("H14").Value = 0.3329182898289247
$ws.Range("I14").Value = 0.9655594880572238
$ws.Range("L14").Value = 0.4463539286308276
$ws.Range("M14").Value = 0.4447482681098478
$ws.Range("O14").Value = 1.038579057994895
$ws.Range("B15").Value = 1.640999551671712
$ws.Range("D15").Value = 0.02219612689257389
$ws.Range("E15").Value = 0.9965685852238124
$ws.Range("F15").Value = 0.3688464131486313
$ws.Range("G15").Value = 0.2292843042613129
$ws.Range("H15").Value = 0.3327707636336328
$ws.Range("I15").Value = 0.9640276430325372
$ws.Range("L15").Value = 0.4421758882129723
$ws.Range("M15").Value = 0.4408977062522723
$ws.Range("O15").Value = 1.035703558925576
$ws.Range("B16").Value = 1.563466841494289
$ws.Range("D16").Value = 0.02082858886578975
$ws.Range("E16").Value = 0.966520397738563
$ws.Range("F16").Value = 0.3609774921561169
$ws.Range("G16").Value = 0.2224983849908
$ws.Range("H16").Value = 0.3319984074341136
$ws.Range("I16").Value = 0.9554681373065961
$ws.Range("L16").Value = 0.4182123295099416
$ws.Range("M16").Value = 0.4188124029195635
$ws.Range("O16").Value = 1.019528136369019
$ws.Range("B17").Value = 1.51584322171459
$ws.Range("D17").Value = 0.01998848088723548
$ws.Range("E17").Value = 0.9481299728548009
$ws.Range("F17").Value = 0.3562354849548512
$ws.Range("G17").Value = 0.2183995268289038
$ws.Range("H17").Value = 0.3315897084185622
$ws.Range("I17").Value = 0.9504128034733839
$ws.Range("L17").Value = 0.4034932995560894
$ws.Range("M17").Value = 0.4052470051638366
$ws.Range("O17").Value = 1.009874552890182
$ws.Range("B18").Value = 1.488429970779578
$ws.Range("D18").Value = 0.01950485655313372
$ws.Range("E18").Value = 0.9375686929602978
$ws.Range("F18").Value = 0.3535392831330739
$ws.Range("G18").Value = 0.2160654079621764
$ws.Range("H18").Value = 0.3313785364104405
$ws.Range("I18").Value = 0.9475771433027376
$ws.Range("L18").Value = 0.3950208211303448
$ws.Range("M18").Value = 0.3974385533940676
$ws.Range("O18").Value = 1.004420806592094
$ws.Range("B19").Value = 1.479144693819251
$ws.Range("D19").Value = 0.01934103934320319
$ws.Range("E19").Value = 0.9339956945270984
$ws.Range("F19").Value = 0.3526317594740433
$ws.Range("G19").Value = 0.2152791296629033
$ws.Range("H19").Value = 0.3313111394081858
$ws.Range("I19").Value = 0.9466294224000649
$ws.Range("L19").Value = 0.3921510906682499
$ws.Range("M19").Value = 0.3947937328118627
$ws.Range("O19").Value = 1.002591191295977
$ws.Range("B20").Value = 1.520915065674956
$ws.Range("D20").Value = 0.02007795520229649
$ws.Range("E20").Value = 0.9500859792887297
$ws.Range("F20").Value = 0.3567370410332629
$ws.Range("G20").Value = 0.21883342915676
$ws.Range("H20").Value = 0.3316307408155126
$ws.Range("I20").Value = 0.9509435001995996
$ws.Range("L20").Value = 0.4050608406219283
$ws.Range("M20").Value = 0.4066916891520904
$ws.Range("O20").Value = 1.010891967234841
$ws.Range("B21").Value = 1.660998013954099
$ws.Range("D21").Value = 0.02254882700466254
$ws.Range("E21").Value = 1.004339340355941
$ws.Range("F21").Value = 0.3709044565123065
$ws.Range("G21").Value = 0.2310562003610386
$ws.Range("H21").Value = 0.3329903259969598
$ws.Range("I21").Value = 0.9662979165003236
$ws.Range("L21").Value = 0.4483570360771409
$ws.Range("M21").Value = 0.4465943698603354
$ws.Range("O21").Value = 1.039963306797176
$ws.Range("B22").Value = 1.752347524212269
$ws.Range("D22").Value = 0.02415970718562477
$ws.Range("E22").Value = 1.039931643986407
$ws.Range("F22").Value = 0.3804438732762065
$ws.Range("G22").Value = 0.2392556076083707
$ws.Range("H22").Value = 0.3340923466224694
$ws.Range("I22").Value = 0.9769715382705613
$ws.Range("L22").Value = 0.4765917983064583
$ws.Range("M22").Value = 0.4726160536161359
$ws.Range("O22").Value = 1.059848851731914
$ws.Range("B23").Value = 1.703611712055704
$ws.Range("D23").Value = 0.02330032703785179
$ws.Range("E23").Value = 1.020923599474685
$ws.Range("F23").Value = 0.3753267883353288
$ws.Range("G23").Value = 0.2348599954184039
$ws.Range("H23").Value = 0.333484616736726
$ws.Range("I23").Value = 0.97121652559224
$ws.Range("L23").Value = 0.4615282258529589
$ws.Range("M23").Value = 0.4587331844108462
$ws.Range("O23").Value = 1.049154204642349
$ws.Range("B24").Value = 1.51862219015635
$ws.Range("D24").Value = 0.02003750583995867
$ws.Range("E24").Value = 0.9492016324535939
$ws.Range("F24").Value = 0.3565101939732713
$ws.Range("G24").Value = 0.2186371922635431
$ws.Range("H24").Value = 0.3316121159567302
$ws.Range("I24").Value = 0.9507033518948873
$ws.Range("L24").Value = 0.4043521873822726
$ws.Range("M24").Value = 0.4060385771458641
$ws.Range("O24").Value = 1.010431694540529
$ws.Range("B25").Value = 1.318269472954171
$ws.Range("D25").Value = 0.01650225177728259
$ws.Range("E25").Value = 0.8724718079843399
$ws.Range("F25").Value = 0.3374080341626851
$ws.Range("G25").Value = 0.2020327200768861
$ws.Range("H25").Value = 0.3305056962784505
$ws.Range("I25").Value = 0.9313322745782884
$ws.Range("L25").Value = 0.3424330323349238
$ws.Range("M25").Value = 0.348971786119229
$ws.Range("O25").Value = 0.9724340452055174
